# Apply weekly update: insert 2 new rows at row 266 (shifting existing rows 266-354 down
# to 268-356), and populate the 2 new rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 266; existing data shifts down.
$ws.Rows.Item(266).Insert()
$ws.Rows.Item(266).Insert()

# New row 266: Zapallo / Camote / 1a nueva(o)
$ws.Cells.Item(266, 1).Value = 10
$ws.Cells.Item(266, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(266, 3).Value = "La Araucanía"
$ws.Cells.Item(266, 4).Value = 44468
$ws.Cells.Item(266, 5).Value = 9
$ws.Cells.Item(266, 6).Value = 100112045
$ws.Cells.Item(266, 7).Value = "Zapallo"
$ws.Cells.Item(266, 8).Value = "Camote"
$ws.Cells.Item(266, 9).Value = "1a nueva(o)"
$ws.Cells.Item(266, 10).Value = 200
$ws.Cells.Item(266, 11).Value = 1000
$ws.Cells.Item(266, 12).Value = 1000
$ws.Cells.Item(266, 13).Value = 1000
$ws.Cells.Item(266, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(266, 15).Value = "Perú"
$ws.Cells.Item(266, 16).Value = 1000
$ws.Cells.Item(266, 17).Value = 1
$ws.Cells.Item(266, 18).Value = "Hortaliza"

# New row 267: Zapallo / Paine / 1a nueva(o)
$ws.Cells.Item(267, 1).Value = 10
$ws.Cells.Item(267, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(267, 3).Value = "La Araucanía"
$ws.Cells.Item(267, 4).Value = 44468
$ws.Cells.Item(267, 5).Value = 9
$ws.Cells.Item(267, 6).Value = 100112045
$ws.Cells.Item(267, 7).Value = "Zapallo"
$ws.Cells.Item(267, 8).Value = "Paine"
$ws.Cells.Item(267, 9).Value = "1a nueva(o)"
$ws.Cells.Item(267, 10).Value = 800
$ws.Cells.Item(267, 11).Value = 500
$ws.Cells.Item(267, 12).Value = 500
$ws.Cells.Item(267, 13).Value = 500
$ws.Cells.Item(267, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(267, 15).Value = "Región del Maule"
$ws.Cells.Item(267, 16).Value = 500
$ws.Cells.Item(267, 17).Value = 1
$ws.Cells.Item(267, 18).Value = "Hortaliza"

# Ensure the date column keeps its date number format (style carried over from Insert,
# but set explicitly to be safe).
$ws.Range("D266:D267").NumberFormat = "YYYY-MM-DD HH:MM:SS"
